$d = $word.ActiveDocument
$notFound = 0

$r = $d.Content
$found = $r.Find.Execute("2023-11-05 Sunday", $true, $false, $false, $false, $false, $false, 1, $false, "2023-11-06 Monday", 1)
if (-not $found) { Write-Output "NOT FOUND: 2023-11-05 Sunday"; $notFound++ }
$r = $d.Content
$found = $r.Find.Execute("58-38=", $true, $false, $false, $false, $false, $false, 1, $false, "4+44=", 1)
if (-not $found) { Write-Output "NOT FOUND: 58-38="; $notFound++ }
$r = $d.Content
$found = $r.Find.Execute("61+18=", $true, $false, $false, $false, $false, $false, 1, $false, "68-2=", 1)
if (-not $found) { Write-Output "NOT FOUND: 61+18="; $notFound++ }
$r = $d.Content
$found = $r.Find.Execute("74-2=", $true, $false, $false, $false, $false, $false, 1, $false, "21+55=", 1)
if (-not $found) { Write-Output "NOT FOUND: 74-2="; $notFound++ }
$r = $d.Content
$found = $r.Find.Execute("51-10=", $true, $false, $false, $false, $false, $false, 1, $false, "17+29=", 1)
if (-not $found) { Write-Output "NOT FOUND: 51-10="; $notFound++ }
$r = $d.Content
$found = $r.Find.Execute("89-69=", $true, $false, $false, $false, $false, $false, 1, $false, "82-0=", 1)
if (-not $found) { Write-Output "NOT FOUND: 89-69="; $notFound++ }
$r = $d.Content
$found = $r.Find.Execute("37+20=", $true, $false, $false, $false, $false, $false, 1, $false, "29+30=", 1)
if (-not $found) { Write-Output "NOT FOUND: 37+20="; $notFound++ }
$r = $d.Content
$found = $r.Find.Execute("8+58=", $true, $false, $false, $false, $false, $false, 1, $false, "19+36=", 1)
if (-not $found) { Write-Output "NOT FOUND: 8+58="; $notFound++ }
$r = $d.Content
$found = $r.Find.Execute("8+21=", $true, $false, $false, $false, $false, $false, 1, $false, "36-8=", 1)
if (-not $found) { Write-Output "NOT FOUND: 8+21="; $notFound++ }
$r = $d.Content
$found = $r.Find.Execute("79-41=", $true, $false, $false, $false, $false, $false, 1, $false, "54-24=", 1)
if (-not $found) { Write-Output "NOT FOUND: 79-41="; $notFound++ }
$r = $d.Content
$found = $r.Find.Execute("13-9=", $true, $false, $false, $false, $false, $false, 1, $false, "42+17=", 1)
if (-not $found) { Write-Output "NOT FOUND: 13-9="; $notFound++ }
$r = $d.Content
$found = $r.Find.Execute("56-52=", $true, $false, $false, $false, $false, $false, 1, $false, "42+1=", 1)
if (-not $found) { Write-Output "NOT FOUND: 56-52="; $notFound++ }
$r = $d.Content
$found = $r.Find.Execute("2+45=", $true, $false, $false, $false, $false, $false, 1, $false, "28+52=", 1)
if (-not $found) { Write-Output "NOT FOUND: 2+45="; $notFound++ }
$r = $d.Content
$found = $r.Find.Execute("48-0=", $true, $false, $false, $false, $false, $false, 1, $false, "13+22=", 1)
if (-not $found) { Write-Output "NOT FOUND: 48-0="; $notFound++ }
$r = $d.Content
$found = $r.Find.Execute("60-27=", $true, $false, $false, $false, $false, $false, 1, $false, "41+3=", 1)
if (-not $found) { Write-Output "NOT FOUND: 60-27="; $notFound++ }
$r = $d.Content
$found = $r.Find.Execute("66-25=", $true, $false, $false, $false, $false, $false, 1, $false, "11+7=", 1)
if (-not $found) { Write-Output "NOT FOUND: 66-25="; $notFound++ }
$r = $d.Content
$found = $r.Find.Execute("97-85=", $true, $false, $false, $false, $false, $false, 1, $false, "91-11=", 1)
if (-not $found) { Write-Output "NOT FOUND: 97-85="; $notFound++ }
$r = $d.Content
$found = $r.Find.Execute("82-12=", $true, $false, $false, $false, $false, $false, 1, $false, "67-62=", 1)
if (-not $found) { Write-Output "NOT FOUND: 82-12="; $notFound++ }
$r = $d.Content
$found = $r.Find.Execute("35+35=", $true, $false, $false, $false, $false, $false, 1, $false, "39-30=", 1)
if (-not $found) { Write-Output "NOT FOUND: 35+35="; $notFound++ }
$r = $d.Content
$found = $r.Find.Execute("37+8=", $true, $false, $false, $false, $false, $false, 1, $false, "85-59=", 1)
if (-not $found) { Write-Output "NOT FOUND: 37+8="; $notFound++ }
$r = $d.Content
$found = $r.Find.Execute("1+36=", $true, $false, $false, $false, $false, $false, 1, $false, "26+29=", 1)
if (-not $found) { Write-Output "NOT FOUND: 1+36="; $notFound++ }
$r = $d.Content
$found = $r.Find.Execute("88-55=", $true, $false, $false, $false, $false, $false, 1, $false, "54-27=", 1)
if (-not $found) { Write-Output "NOT FOUND: 88-55="; $notFound++ }
$r = $d.Content
$found = $r.Find.Execute("43-26=", $true, $false, $false, $false, $false, $false, 1, $false, "47+18=", 1)
if (-not $found) { Write-Output "NOT FOUND: 43-26="; $notFound++ }
$r = $d.Content
$found = $r.Find.Execute("64-7=", $true, $false, $false, $false, $false, $false, 1, $false, "24+28=", 1)
if (-not $found) { Write-Output "NOT FOUND: 64-7="; $notFound++ }
$r = $d.Content
$found = $r.Find.Execute("2+96=", $true, $false, $false, $false, $false, $false, 1, $false, "28+51=", 1)
if (-not $found) { Write-Output "NOT FOUND: 2+96="; $notFound++ }
$r = $d.Content
$found = $r.Find.Execute("13+81=", $true, $false, $false, $false, $false, $false, 1, $false, "81-56=", 1)
if (-not $found) { Write-Output "NOT FOUND: 13+81="; $notFound++ }
$r = $d.Content
$found = $r.Find.Execute("20-14=", $true, $false, $false, $false, $false, $false, 1, $false, "72-36=", 1)
if (-not $found) { Write-Output "NOT FOUND: 20-14="; $notFound++ }
$r = $d.Content
$found = $r.Find.Execute("81+15=", $true, $false, $false, $false, $false, $false, 1, $false, "84-8=", 1)
if (-not $found) { Write-Output "NOT FOUND: 81+15="; $notFound++ }
$r = $d.Content
$found = $r.Find.Execute("27+15=", $true, $false, $false, $false, $false, $false, 1, $false, "17+28=", 1)
if (-not $found) { Write-Output "NOT FOUND: 27+15="; $notFound++ }
$r = $d.Content
$found = $r.Find.Execute("81-69=", $true, $false, $false, $false, $false, $false, 1, $false, "16+56=", 1)
if (-not $found) { Write-Output "NOT FOUND: 81-69="; $notFound++ }
$r = $d.Content
$found = $r.Find.Execute("84-3=", $true, $false, $false, $false, $false, $false, 1, $false, "67-0=", 1)
if (-not $found) { Write-Output "NOT FOUND: 84-3="; $notFound++ }
$r = $d.Content
$found = $r.Find.Execute("59-22=", $true, $false, $false, $false, $false, $false, 1, $false, "41+2=", 1)
if (-not $found) { Write-Output "NOT FOUND: 59-22="; $notFound++ }
$r = $d.Content
$found = $r.Find.Execute("14+26=", $true, $false, $false, $false, $false, $false, 1, $false, "83-43=", 1)
if (-not $found) { Write-Output "NOT FOUND: 14+26="; $notFound++ }
$r = $d.Content
$found = $r.Find.Execute("57+18=", $true, $false, $false, $false, $false, $false, 1, $false, "96-57=", 1)
if (-not $found) { Write-Output "NOT FOUND: 57+18="; $notFound++ }
$r = $d.Content
$found = $r.Find.Execute("17+57=", $true, $false, $false, $false, $false, $false, 1, $false, "12-3=", 1)
if (-not $found) { Write-Output "NOT FOUND: 17+57="; $notFound++ }
$r = $d.Content
$found = $r.Find.Execute("12+15=", $true, $false, $false, $false, $false, $false, 1, $false, "97-50=", 1)
if (-not $found) { Write-Output "NOT FOUND: 12+15="; $notFound++ }
$r = $d.Content
$found = $r.Find.Execute("14+26=", $true, $false, $false, $false, $false, $false, 1, $false, "1+46=", 1)
if (-not $found) { Write-Output "NOT FOUND: 14+26="; $notFound++ }
$r = $d.Content
$found = $r.Find.Execute("14+12=", $true, $false, $false, $false, $false, $false, 1, $false, "73-39=", 1)
if (-not $found) { Write-Output "NOT FOUND: 14+12="; $notFound++ }
$r = $d.Content
$found = $r.Find.Execute("61+6=", $true, $false, $false, $false, $false, $false, 1, $false, "41+39=", 1)
if (-not $found) { Write-Output "NOT FOUND: 61+6="; $notFound++ }
$r = $d.Content
$found = $r.Find.Execute("80-2=", $true, $false, $false, $false, $false, $false, 1, $false, "97-47=", 1)
if (-not $found) { Write-Output "NOT FOUND: 80-2="; $notFound++ }
$r = $d.Content
$found = $r.Find.Execute("25-9=", $true, $false, $false, $false, $false, $false, 1, $false, "92-69=", 1)
if (-not $found) { Write-Output "NOT FOUND: 25-9="; $notFound++ }
$r = $d.Content
$found = $r.Find.Execute("27+6=", $true, $false, $false, $false, $false, $false, 1, $false, "2+9=", 1)
if (-not $found) { Write-Output "NOT FOUND: 27+6="; $notFound++ }
$r = $d.Content
$found = $r.Find.Execute("33+49=", $true, $false, $false, $false, $false, $false, 1, $false, "7+19=", 1)
if (-not $found) { Write-Output "NOT FOUND: 33+49="; $notFound++ }
$r = $d.Content
$found = $r.Find.Execute("66-58=", $true, $false, $false, $false, $false, $false, 1, $false, "64-1=", 1)
if (-not $found) { Write-Output "NOT FOUND: 66-58="; $notFound++ }
$r = $d.Content
$found = $r.Find.Execute("41+37=", $true, $false, $false, $false, $false, $false, 1, $false, "9+47=", 1)
if (-not $found) { Write-Output "NOT FOUND: 41+37="; $notFound++ }
$r = $d.Content
$found = $r.Find.Execute("32+46=", $true, $false, $false, $false, $false, $false, 1, $false, "1+72=", 1)
if (-not $found) { Write-Output "NOT FOUND: 32+46="; $notFound++ }
$r = $d.Content
$found = $r.Find.Execute("3+86=", $true, $false, $false, $false, $false, $false, 1, $false, "8+20=", 1)
if (-not $found) { Write-Output "NOT FOUND: 3+86="; $notFound++ }
$r = $d.Content
$found = $r.Find.Execute("67-63=", $true, $false, $false, $false, $false, $false, 1, $false, "73-66=", 1)
if (-not $found) { Write-Output "NOT FOUND: 67-63="; $notFound++ }
$r = $d.Content
$found = $r.Find.Execute("13+84=", $true, $false, $false, $false, $false, $false, 1, $false, "79-74=", 1)
if (-not $found) { Write-Output "NOT FOUND: 13+84="; $notFound++ }
$r = $d.Content
$found = $r.Find.Execute("15+35=", $true, $false, $false, $false, $false, $false, 1, $false, "17+7=", 1)
if (-not $found) { Write-Output "NOT FOUND: 15+35="; $notFound++ }
$r = $d.Content
$found = $r.Find.Execute("76-6=", $true, $false, $false, $false, $false, $false, 1, $false, "4+68=", 1)
if (-not $found) { Write-Output "NOT FOUND: 76-6="; $notFound++ }
$r = $d.Content
$found = $r.Find.Execute("29+17=", $true, $false, $false, $false, $false, $false, 1, $false, "67-67=", 1)
if (-not $found) { Write-Output "NOT FOUND: 29+17="; $notFound++ }
$r = $d.Content
$found = $r.Find.Execute("0+94=", $true, $false, $false, $false, $false, $false, 1, $false, "4+57=", 1)
if (-not $found) { Write-Output "NOT FOUND: 0+94="; $notFound++ }
$r = $d.Content
$found = $r.Find.Execute("25+55=", $true, $false, $false, $false, $false, $false, 1, $false, "79-9=", 1)
if (-not $found) { Write-Output "NOT FOUND: 25+55="; $notFound++ }
$r = $d.Content
$found = $r.Find.Execute("31-18=", $true, $false, $false, $false, $false, $false, 1, $false, "64+2=", 1)
if (-not $found) { Write-Output "NOT FOUND: 31-18="; $notFound++ }
$r = $d.Content
$found = $r.Find.Execute("0+86=", $true, $false, $false, $false, $false, $false, 1, $false, "14+83=", 1)
if (-not $found) { Write-Output "NOT FOUND: 0+86="; $notFound++ }
$r = $d.Content
$found = $r.Find.Execute("70-49=", $true, $false, $false, $false, $false, $false, 1, $false, "62-50=", 1)
if (-not $found) { Write-Output "NOT FOUND: 70-49="; $notFound++ }
$r = $d.Content
$found = $r.Find.Execute("3+86=", $true, $false, $false, $false, $false, $false, 1, $false, "82-27=", 1)
if (-not $found) { Write-Output "NOT FOUND: 3+86="; $notFound++ }
$r = $d.Content
$found = $r.Find.Execute("47+41=", $true, $false, $false, $false, $false, $false, 1, $false, "52+46=", 1)
if (-not $found) { Write-Output "NOT FOUND: 47+41="; $notFound++ }
$r = $d.Content
$found = $r.Find.Execute("58+36=", $true, $false, $false, $false, $false, $false, 1, $false, "38-37=", 1)
if (-not $found) { Write-Output "NOT FOUND: 58+36="; $notFound++ }
$r = $d.Content
$found = $r.Find.Execute("32-27=", $true, $false, $false, $false, $false, $false, 1, $false, "88-52=", 1)
if (-not $found) { Write-Output "NOT FOUND: 32-27="; $notFound++ }
$r = $d.Content
$found = $r.Find.Execute("80+15=", $true, $false, $false, $false, $false, $false, 1, $false, "31-5=", 1)
if (-not $found) { Write-Output "NOT FOUND: 80+15="; $notFound++ }
$r = $d.Content
$found = $r.Find.Execute("3+91=", $true, $false, $false, $false, $false, $false, 1, $false, "67+31=", 1)
if (-not $found) { Write-Output "NOT FOUND: 3+91="; $notFound++ }
$r = $d.Content
$found = $r.Find.Execute("30+11=", $true, $false, $false, $false, $false, $false, 1, $false, "67-49=", 1)
if (-not $found) { Write-Output "NOT FOUND: 30+11="; $notFound++ }
$r = $d.Content
$found = $r.Find.Execute("13+28=", $true, $false, $false, $false, $false, $false, 1, $false, "14+2=", 1)
if (-not $found) { Write-Output "NOT FOUND: 13+28="; $notFound++ }
$r = $d.Content
$found = $r.Find.Execute("62-0=", $true, $false, $false, $false, $false, $false, 1, $false, "75-64=", 1)
if (-not $found) { Write-Output "NOT FOUND: 62-0="; $notFound++ }
$r = $d.Content
$found = $r.Find.Execute("47-6=", $true, $false, $false, $false, $false, $false, 1, $false, "64+13=", 1)
if (-not $found) { Write-Output "NOT FOUND: 47-6="; $notFound++ }
$r = $d.Content
$found = $r.Find.Execute("66-36=", $true, $false, $false, $false, $false, $false, 1, $false, "22+46=", 1)
if (-not $found) { Write-Output "NOT FOUND: 66-36="; $notFound++ }
$r = $d.Content
$found = $r.Find.Execute("43+12=", $true, $false, $false, $false, $false, $false, 1, $false, "49+1=", 1)
if (-not $found) { Write-Output "NOT FOUND: 43+12="; $notFound++ }
$r = $d.Content
$found = $r.Find.Execute("5+21=", $true, $false, $false, $false, $false, $false, 1, $false, "97-46=", 1)
if (-not $found) { Write-Output "NOT FOUND: 5+21="; $notFound++ }
$r = $d.Content
$found = $r.Find.Execute("63-39=", $true, $false, $false, $false, $false, $false, 1, $false, "27-8=", 1)
if (-not $found) { Write-Output "NOT FOUND: 63-39="; $notFound++ }
$r = $d.Content
$found = $r.Find.Execute("28+47=", $true, $false, $false, $false, $false, $false, 1, $false, "27+40=", 1)
if (-not $found) { Write-Output "NOT FOUND: 28+47="; $notFound++ }
$r = $d.Content
$found = $r.Find.Execute("60-10=", $true, $false, $false, $false, $false, $false, 1, $false, "21+21=", 1)
if (-not $found) { Write-Output "NOT FOUND: 60-10="; $notFound++ }
$r = $d.Content
$found = $r.Find.Execute("24-7=", $true, $false, $false, $false, $false, $false, 1, $false, "61-7=", 1)
if (-not $found) { Write-Output "NOT FOUND: 24-7="; $notFound++ }
$r = $d.Content
$found = $r.Find.Execute("55+33=", $true, $false, $false, $false, $false, $false, 1, $false, "86-75=", 1)
if (-not $found) { Write-Output "NOT FOUND: 55+33="; $notFound++ }
$r = $d.Content
$found = $r.Find.Execute("73-23=", $true, $false, $false, $false, $false, $false, 1, $false, "79-35=", 1)
if (-not $found) { Write-Output "NOT FOUND: 73-23="; $notFound++ }
$r = $d.Content
$found = $r.Find.Execute("32+51=", $true, $false, $false, $false, $false, $false, 1, $false, "16+40=", 1)
if (-not $found) { Write-Output "NOT FOUND: 32+51="; $notFound++ }
$r = $d.Content
$found = $r.Find.Execute("61-3=", $true, $false, $false, $false, $false, $false, 1, $false, "96-69=", 1)
if (-not $found) { Write-Output "NOT FOUND: 61-3="; $notFound++ }
$r = $d.Content
$found = $r.Find.Execute("3-1=", $true, $false, $false, $false, $false, $false, 1, $false, "55-25=", 1)
if (-not $found) { Write-Output "NOT FOUND: 3-1="; $notFound++ }
$r = $d.Content
$found = $r.Find.Execute("47+36=", $true, $false, $false, $false, $false, $false, 1, $false, "82-64=", 1)
if (-not $found) { Write-Output "NOT FOUND: 47+36="; $notFound++ }
$r = $d.Content
$found = $r.Find.Execute("24+50=", $true, $false, $false, $false, $false, $false, 1, $false, "77-12=", 1)
if (-not $found) { Write-Output "NOT FOUND: 24+50="; $notFound++ }
$r = $d.Content
$found = $r.Find.Execute("24+56=", $true, $false, $false, $false, $false, $false, 1, $false, "89-58=", 1)
if (-not $found) { Write-Output "NOT FOUND: 24+56="; $notFound++ }
$r = $d.Content
$found = $r.Find.Execute("45+2=", $true, $false, $false, $false, $false, $false, 1, $false, "20+25=", 1)
if (-not $found) { Write-Output "NOT FOUND: 45+2="; $notFound++ }
$r = $d.Content
$found = $r.Find.Execute("17+34=", $true, $false, $false, $false, $false, $false, 1, $false, "85-85=", 1)
if (-not $found) { Write-Output "NOT FOUND: 17+34="; $notFound++ }
$r = $d.Content
$found = $r.Find.Execute("78-67=", $true, $false, $false, $false, $false, $false, 1, $false, "52-1=", 1)
if (-not $found) { Write-Output "NOT FOUND: 78-67="; $notFound++ }
$r = $d.Content
$found = $r.Find.Execute("95-91=", $true, $false, $false, $false, $false, $false, 1, $false, "47-1=", 1)
if (-not $found) { Write-Output "NOT FOUND: 95-91="; $notFound++ }
$r = $d.Content
$found = $r.Find.Execute("86-66=", $true, $false, $false, $false, $false, $false, 1, $false, "69-13=", 1)
if (-not $found) { Write-Output "NOT FOUND: 86-66="; $notFound++ }
$r = $d.Content
$found = $r.Find.Execute("97-95=", $true, $false, $false, $false, $false, $false, 1, $false, "9+19=", 1)
if (-not $found) { Write-Output "NOT FOUND: 97-95="; $notFound++ }
$r = $d.Content
$found = $r.Find.Execute("64+29=", $true, $false, $false, $false, $false, $false, 1, $false, "47-37=", 1)
if (-not $found) { Write-Output "NOT FOUND: 64+29="; $notFound++ }
$r = $d.Content
$found = $r.Find.Execute("98-0=", $true, $false, $false, $false, $false, $false, 1, $false, "94-4=", 1)
if (-not $found) { Write-Output "NOT FOUND: 98-0="; $notFound++ }
$r = $d.Content
$found = $r.Find.Execute("55+29=", $true, $false, $false, $false, $false, $false, 1, $false, "56-50=", 1)
if (-not $found) { Write-Output "NOT FOUND: 55+29="; $notFound++ }
$r = $d.Content
$found = $r.Find.Execute("45-27=", $true, $false, $false, $false, $false, $false, 1, $false, "74+12=", 1)
if (-not $found) { Write-Output "NOT FOUND: 45-27="; $notFound++ }
$r = $d.Content
$found = $r.Find.Execute("37+62=", $true, $false, $false, $false, $false, $false, 1, $false, "17+19=", 1)
if (-not $found) { Write-Output "NOT FOUND: 37+62="; $notFound++ }
$r = $d.Content
$found = $r.Find.Execute("85+4=", $true, $false, $false, $false, $false, $false, 1, $false, "66+3=", 1)
if (-not $found) { Write-Output "NOT FOUND: 85+4="; $notFound++ }
$r = $d.Content
$found = $r.Find.Execute("28+53=", $true, $false, $false, $false, $false, $false, 1, $false, "29+4=", 1)
if (-not $found) { Write-Output "NOT FOUND: 28+53="; $notFound++ }
$r = $d.Content
$found = $r.Find.Execute("77-14=", $true, $false, $false, $false, $false, $false, 1, $false, "89-28=", 1)
if (-not $found) { Write-Output "NOT FOUND: 77-14="; $notFound++ }
$r = $d.Content
$found = $r.Find.Execute("15+60=", $true, $false, $false, $false, $false, $false, 1, $false, "97-33=", 1)
if (-not $found) { Write-Output "NOT FOUND: 15+60="; $notFound++ }
$r = $d.Content
$found = $r.Find.Execute("3+47=", $true, $false, $false, $false, $false, $false, 1, $false, "61+19=", 1)
if (-not $found) { Write-Output "NOT FOUND: 3+47="; $notFound++ }
$r = $d.Content
$found = $r.Find.Execute("85-16=", $true, $false, $false, $false, $false, $false, 1, $false, "47+38=", 1)
if (-not $found) { Write-Output "NOT FOUND: 85-16="; $notFound++ }
$r = $d.Content
$found = $r.Find.Execute("83+8=", $true, $false, $false, $false, $false, $false, 1, $false, "99-1=", 1)
if (-not $found) { Write-Output "NOT FOUND: 83+8="; $notFound++ }
$r = $d.Content
$found = $r.Find.Execute("21+53=", $true, $false, $false, $false, $false, $false, 1, $false, "96-11=", 1)
if (-not $found) { Write-Output "NOT FOUND: 21+53="; $notFound++ }
Write-Output "done, notFound=$notFound"